$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("10796", "Nom non trouvé", "165", "N/A", "https://www.lego.com/cdn/cs/set/assets/blt8ceb255a3e66c253/10796.png?format=webply&fit=bounds&quality=75&width=528&height=528&dpr=1", "https://www.lego.com/fr-fr/product/10796"),
    @("77246", "Nom non trouvé", "248", "N/A", "https://www.lego.com/cdn/cs/set/assets/blt0a45cbfbac56a36b/77246_Prod_en-gb.png?format=webply&fit=bounds&quality=75&width=528&height=528&dpr=1", "https://www.lego.com/fr-fr/product/77246"),
    @("42154", "Nom non trouvé", "1468", "N/A", "https://www.lego.com/cdn/cs/set/assets/blt5014244d8d8dc8ad/42154.png?format=webply&fit=bounds&quality=75&width=528&height=528&dpr=1", "https://www.lego.com/fr-fr/product/42154")
)

$startRow = 27
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    # Prefix numeric-looking text (ID_Set, nbPieces) with an apostrophe so
    # Excel stores it as text rather than auto-converting to a number -
    # this matches every other row in the sheet, where these columns are
    # text, not numeric, cells.
    $ws.Cells.Item($row, 1).Value = "'" + $rowData[0]
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = "'" + $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
    $ws.Cells.Item($row, 5).Value = $rowData[4]
    $ws.Cells.Item($row, 6).Value = $rowData[5]
    # Remaining retailer-URL columns have no data yet for these sets - keep
    # them as empty text cells (matching every other row's untouched
    # trailing columns).
    $ws.Cells.Item($row, 7).Value = "'"
    $ws.Cells.Item($row, 8).Value = "'"
    $ws.Cells.Item($row, 9).Value = "'"
    $ws.Cells.Item($row, 10).Value = "'"

    # The leading apostrophe above forces text storage for the numeric-
    # looking / empty cells, but also stamps them with a "quote prefix"
    # style. Re-apply the plain Normal style so the new cells carry no
    # explicit formatting, same as the rest of the data rows.
    $rowRange = $ws.Range($ws.Cells.Item($row, 1), $ws.Cells.Item($row, 10))
    $rowRange.Style = "Normal"
}
